# Natmi following Dr Hou advice
# Re-run of the Wnt5a-Fzd3 LR-pair NATMI export: the sending/target cluster
# set now includes the "sCs" and "ECs" populations in addition to "FAPs",
# so the ligand-receptor edge table is recomputed and grows from 2 data
# rows to 6 data rows (A2:T7), with refreshed statistics on existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: FAPs -> Wnt5a -> Fzd3 -> FAPs
$ws.Cells.Item(2, "A").Value = "FAPs"
$ws.Cells.Item(2, "B").Value = "Wnt5a"
$ws.Cells.Item(2, "C").Value = "Fzd3"
$ws.Cells.Item(2, "D").Value = "FAPs"
$ws.Cells.Item(2, "E").Value = 3
$ws.Cells.Item(2, "F").Value = 1
$ws.Cells.Item(2, "G").Value = 10.43424333333333
$ws.Cells.Item(2, "H").Value = 31.30273
$ws.Cells.Item(2, "I").Value = 0.9711091978791583
$ws.Cells.Item(2, "J").Value = 0.9711091978791584
$ws.Cells.Item(2, "K").Value = 2
$ws.Cells.Item(2, "L").Value = 0.6666666666666666
$ws.Cells.Item(2, "M").Value = 0.173174
$ws.Cells.Item(2, "N").Value = 0.519522
$ws.Cells.Item(2, "O").Value = 0.0473074116693291
$ws.Cells.Item(2, "P").Value = 0.0473074116693291
$ws.Cells.Item(2, "Q").Value = 1.806939655006667
$ws.Cells.Item(2, "R").Value = 16.26245689506
$ws.Cells.Item(2, "S").Value = 0.04594066259994131
$ws.Cells.Item(2, "T").Value = 0.04594066259994132

# Row 3: FAPs -> Wnt5a -> Fzd3 -> sCs
$ws.Cells.Item(3, "A").Value = "FAPs"
$ws.Cells.Item(3, "B").Value = "Wnt5a"
$ws.Cells.Item(3, "C").Value = "Fzd3"
$ws.Cells.Item(3, "D").Value = "sCs"
$ws.Cells.Item(3, "E").Value = 3
$ws.Cells.Item(3, "F").Value = 1
$ws.Cells.Item(3, "G").Value = 10.43424333333333
$ws.Cells.Item(3, "H").Value = 31.30273
$ws.Cells.Item(3, "I").Value = 0.9711091978791583
$ws.Cells.Item(3, "J").Value = 0.9711091978791584
$ws.Cells.Item(3, "K").Value = 3
$ws.Cells.Item(3, "L").Value = 1
$ws.Cells.Item(3, "M").Value = 0.880936
$ws.Cells.Item(3, "N").Value = 2.642808
$ws.Cells.Item(3, "O").Value = 0.2406527654632456
$ws.Cells.Item(3, "P").Value = 0.2406527654632456
$ws.Cells.Item(3, "Q").Value = 9.191900585093334
$ws.Cells.Item(3, "R").Value = 82.72710526584
$ws.Cells.Item(3, "S").Value = 0.2337001140364136
$ws.Cells.Item(3, "T").Value = 0.2337001140364137

# Row 4 (new): FAPs -> Wnt5a -> Fzd3 -> ECs
$ws.Cells.Item(4, "A").Value = "FAPs"
$ws.Cells.Item(4, "B").Value = "Wnt5a"
$ws.Cells.Item(4, "C").Value = "Fzd3"
$ws.Cells.Item(4, "D").Value = "ECs"
$ws.Cells.Item(4, "E").Value = 3
$ws.Cells.Item(4, "F").Value = 1
$ws.Cells.Item(4, "G").Value = 10.43424333333333
$ws.Cells.Item(4, "H").Value = 31.30273
$ws.Cells.Item(4, "I").Value = 0.9711091978791583
$ws.Cells.Item(4, "J").Value = 0.9711091978791584
$ws.Cells.Item(4, "K").Value = 3
$ws.Cells.Item(4, "L").Value = 1
$ws.Cells.Item(4, "M").Value = 2.606500333333333
$ws.Cells.Item(4, "N").Value = 7.819501
$ws.Cells.Item(4, "O").Value = 0.7120398228674253
$ws.Cells.Item(4, "P").Value = 0.7120398228674253
$ws.Cells.Item(4, "Q").Value = 27.19685872641444
$ws.Cells.Item(4, "R").Value = 244.77172853773
$ws.Cells.Item(4, "S").Value = 0.6914684212428033
$ws.Cells.Item(4, "T").Value = 0.6914684212428034

# Row 5 (new): sCs -> Wnt5a -> Fzd3 -> ECs
$ws.Cells.Item(5, "A").Value = "sCs"
$ws.Cells.Item(5, "B").Value = "Wnt5a"
$ws.Cells.Item(5, "C").Value = "Fzd3"
$ws.Cells.Item(5, "D").Value = "ECs"
$ws.Cells.Item(5, "E").Value = 2
$ws.Cells.Item(5, "F").Value = 0.6666666666666666
$ws.Cells.Item(5, "G").Value = 0.310422
$ws.Cells.Item(5, "H").Value = 0.9312659999999999
$ws.Cells.Item(5, "I").Value = 0.02889080212084161
$ws.Cells.Item(5, "J").Value = 0.02889080212084161
$ws.Cells.Item(5, "K").Value = 2
$ws.Cells.Item(5, "L").Value = 0.6666666666666666
$ws.Cells.Item(5, "M").Value = 0.173174
$ws.Cells.Item(5, "N").Value = 0.519522
$ws.Cells.Item(5, "O").Value = 0.0473074116693291
$ws.Cells.Item(5, "P").Value = 0.0473074116693291
$ws.Cells.Item(5, "Q").Value = 0.053757019428
$ws.Cells.Item(5, "R").Value = 0.483813174852
$ws.Cells.Item(5, "S").Value = 0.00136674906938778
$ws.Cells.Item(5, "T").Value = 0.00136674906938778

# Row 6 (new): sCs -> Wnt5a -> Fzd3 -> FAPs
$ws.Cells.Item(6, "A").Value = "sCs"
$ws.Cells.Item(6, "B").Value = "Wnt5a"
$ws.Cells.Item(6, "C").Value = "Fzd3"
$ws.Cells.Item(6, "D").Value = "FAPs"
$ws.Cells.Item(6, "E").Value = 2
$ws.Cells.Item(6, "F").Value = 0.6666666666666666
$ws.Cells.Item(6, "G").Value = 0.310422
$ws.Cells.Item(6, "H").Value = 0.9312659999999999
$ws.Cells.Item(6, "I").Value = 0.02889080212084161
$ws.Cells.Item(6, "J").Value = 0.02889080212084161
$ws.Cells.Item(6, "K").Value = 3
$ws.Cells.Item(6, "L").Value = 1
$ws.Cells.Item(6, "M").Value = 0.880936
$ws.Cells.Item(6, "N").Value = 2.642808
$ws.Cells.Item(6, "O").Value = 0.2406527654632456
$ws.Cells.Item(6, "P").Value = 0.2406527654632456
$ws.Cells.Item(6, "Q").Value = 0.273461914992
$ws.Cells.Item(6, "R").Value = 2.461157234928
$ws.Cells.Item(6, "S").Value = 0.006952651426831934
$ws.Cells.Item(6, "T").Value = 0.006952651426831936

# Row 7 (new): sCs -> Wnt5a -> Fzd3 -> sCs
$ws.Cells.Item(7, "A").Value = "sCs"
$ws.Cells.Item(7, "B").Value = "Wnt5a"
$ws.Cells.Item(7, "C").Value = "Fzd3"
$ws.Cells.Item(7, "D").Value = "sCs"
$ws.Cells.Item(7, "E").Value = 2
$ws.Cells.Item(7, "F").Value = 0.6666666666666666
$ws.Cells.Item(7, "G").Value = 0.310422
$ws.Cells.Item(7, "H").Value = 0.9312659999999999
$ws.Cells.Item(7, "I").Value = 0.02889080212084161
$ws.Cells.Item(7, "J").Value = 0.02889080212084161
$ws.Cells.Item(7, "K").Value = 3
$ws.Cells.Item(7, "L").Value = 1
$ws.Cells.Item(7, "M").Value = 2.606500333333333
$ws.Cells.Item(7, "N").Value = 7.819501
$ws.Cells.Item(7, "O").Value = 0.7120398228674253
$ws.Cells.Item(7, "P").Value = 0.7120398228674253
$ws.Cells.Item(7, "Q").Value = 0.809115046474
$ws.Cells.Item(7, "R").Value = 7.282035418265999
$ws.Cells.Item(7, "S").Value = 0.02057140162462189
$ws.Cells.Item(7, "T").Value = 0.0205714016246219
